# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Shift the existing MSME table (rows 5-8) down to rows 11, 13, 14, 15,
# leaving gaps for the new "Source Type" label (row 9) and the new
# "Employment (% of total)" data row (row 12).
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(12).Insert()

# Insert() copies formatting from the row above into the blank row; drop
# the phantom styled-but-empty cells it created at B12:D12.
$ws.Range("B12:D12").Clear()

# New "Source Type" label above the first table (bold + underline)
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# New "Employment (% of total)" row inside the first table
$ws.Range("A12").Value = "Employment (% of total)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32"
$ws.Range("D12").ClearFormats()

# Second table: headers
$ws.Range("B17").Value = "Micro"
$ws.Range("B17").Font.Bold = $true
$ws.Range("C17").Value = "SMEs"
$ws.Range("C17").Font.Bold = $true
$ws.Range("D17").Value = "MSMEs"
$ws.Range("D17").Font.Bold = $true

# Second table: data row
$ws.Range("A18").Value = "Value added to the economy (% of total)"
$ws.Range("A18").Font.Bold = $true
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45"
$ws.Range("D18").ClearFormats()

# Second table: source footnote (italic, matching the first table's source style)
$ws.Range("A19").Value = "Source: BELTRAIDE, 2001"
$ws.Range("A19").Font.Italic = $true
